$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set new value in C5 (new shared string "final_part1_log1")
$ws.Range("C5").Value = "final_part1_log1"

# Update the selected cell/range to C6 (as in the diff)
$ws.Range("C6").Select()
